# "Generate Report for Handback"
#
# The file 6c476987-2e90-41aa-b3d3-52976c12f200.md has been handed back for
# both the zh-cn and de-de locales: its status moves from "Ready for
# handoff" to "Handed back: in sync with en-US", and a new "Latest Handback
# DateTime" is recorded on each locale's worksheet.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet: row for 6c476987-...md (row 3) ---
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

# --- zh-cn sheet: row for 6c476987-...md (row 3) ---
$wsZhCn.Range("B3").Value = $newStatus
$wsZhCn.Range("G3").Value = "2016-03-09 18:42:10"

# --- de-de sheet: row for 6c476987-...md (row 3) ---
$wsDeDe.Range("B3").Value = $newStatus
$wsDeDe.Range("G3").Value = "2016-03-09 18:42:22"
